$wb = $excel.ActiveWorkbook

# --- Summary sheet: Generated-on timestamp + Dictionary of Results ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B9").Value2 = "2025-01-08 12:57:37"
$summary.Range("B20").Value2 = '{''mcar'': {''mean'': {''bayesian'': {''0.05'': {''scores'': {''RMSE'': 0.5197922283008971, ''MAE'': 0.4543356516868202, ''MI'': 0.0, ''CORRELATION'': 0}, ''times'': {''contamination'': 0.0004885196685791016, ''optimization'': 0, ''imputation'': 0.0002512931823730469, ''log_imputation'': -8.28889024321749}}, ''0.1'': {''scores'': {''RMSE'': 1.0659202645786816, ''MAE'': 0.9085417731383956, ''MI'': 0.0, ''CORRELATION'': 0}, ''times'': {''contamination'': 0.0003097057342529297, ''optimization'': 0, ''imputation'': 0.00014352798461914062, ''log_imputation'': -8.848980527009976}}, ''0.2'': {''scores'': {''RMSE'': 1.1400385999631493, ''MAE'': 0.9394950730289477, ''MI'': 0.0, ''CORRELATION'': 0}, ''times'': {''contamination'': 0.00023746490478515625, ''optimization'': 0, ''imputation'': 9.608268737792969e-05, ''log_imputation'': -9.250301410372114}}, ''0.4'': {''scores'': {''RMSE'': 1.0333061850175014, ''MAE'': 0.817160720129779, ''MI'': 0.0, ''CORRELATION'': 0}, ''times'': {''contamination'': 0.00040912628173828125, ''optimization'': 0, ''imputation'': 7.462501525878906e-05, ''log_imputation'': -9.503034781778643}}, ''0.6'': {''scores'': {''RMSE'': 1.0938413270459857, ''MAE'': 0.8545290213993658, ''MI'': 0.0, ''CORRELATION'': 0}, ''times'': {''contamination'': 0.0006306171417236328, ''optimization'': 0, ''imputation'': 9.202957153320312e-05, ''log_imputation'': -9.293400602853966}}, ''0.8'': {''scores'': {''RMSE'': 1.07436956341757, ''MAE'': 0.8291370178635111, ''MI'': 0.0, ''CORRELATION'': 0}, ''times'': {''contamination'': 0.001354217529296875, ''optimization'': 0, ''imputation'': 6.961822509765625e-05, ''log_imputation'': -9.572484170050515}}}}, ''cdrec'': {''bayesian'': {''0.05'': {''scores'': {''RMSE'': 0.37483452324301586, ''MAE'': 0.3375262694281006, ''MI'': 1.0397207708399179, ''CORRELATION'': 0.7365655689896633}, ''times'': {''contamination'': 0.00011134147644042969, ''optimization'': 0.5079936981201172, ''imputation'': 0.004687070846557617, ''log_imputation'': -5.362947444611191}}, ''0.1'': {''scores'': {''RMSE'': 1.3799678230195285, ''MAE'': 1.1003322284844623, ''MI'': 1.732867951399863, ''CORRELATION'': -0.500100644242659}, ''times'': {''contamination'': 0.0003001689910888672, ''optimization'': 0.5079936981201172, ''imputation'': 0.0036382675170898438, ''log_imputation'': -5.616247667480349}}, ''0.2'': {''scores'': {''RMSE'': 0.5279485898506157, ''MAE'': 0.42431581904234256, ''MI'': 1.342409426595628, ''CORRELATION'': 0.9071070625126642}, ''times'': {''contamination'': 0.00042891502380371094, ''optimization'': 0.5079936981201172, ''imputation'': 0.015071868896484375, ''log_imputation'': -4.194925259667625}}, ''0.4'': {''scores'': {''RMSE'': 0.6529812630837011, ''MAE'': 0.42858056477338186, ''MI'': 0.8905639332827393, ''CORRELATION'': 0.7722811146383882}, ''times'': {''contamination'': 0.00032210350036621094, ''optimization'': 0.5079936981201172, ''imputation'': 0.008674860000610352, ''log_imputation'': -4.74732609153761}}, ''0.6'': {''scores'': {''RMSE'': 0.6798826891423311, ''MAE'': 0.47135122101632737, ''MI'': 0.6001676421795947, ''CORRELATION'': 0.7742382236368857}, ''times'': {''contamination'': 0.0015621185302734375, ''optimization'': 0.5079936981201172, ''imputation'': 0.009982585906982422, ''log_imputation'': -4.606913113305608}}, ''0.8'': {''scores'': {''RMSE'': 0.7608485588056992, ''MAE'': 0.5479154581689161, ''MI'': 0.42721564894947844, ''CORRELATION'': 0.7017141157422242}, ''times'': {''contamination'': 0.003453493118286133, ''optimization'': 0.5079936981201172, ''imputation'': 0.0073740482330322266, ''log_imputation'': -4.9097884382325}}}}, ''stmvl'': {''bayesian'': {''0.05'': {''scores'': {''RMSE'': 0.3251125774837754, ''MAE'': 0.26797673641099284, ''MI'': 1.0397207708399179, ''CORRELATION'': 0.6142581896031455}, ''times'': {''contamination'': 0.00034546852111816406, ''optimization'': 0.23014473915100098, ''imputation'': 0.06188774108886719, ''log_imputation'': -2.7824331627034407}}, ''0.1'': {''scores'': {''RMSE'': 0.299492451492057, ''MAE'': 0.26432871720074347, ''MI'': 1.9061547465398494, ''CORRELATION'': 0.967896575643492}, ''times'': {''contamination'': 0.0001323223114013672, ''optimization'': 0.23014473915100098, ''imputation'': 0.049030303955078125, ''log_imputation'': -3.015316723968884}}, ''0.2'': {''scores'': {''RMSE'': 0.32852543256899075, ''MAE'': 0.27202573018354975, ''MI'': 1.5996631161656454, ''CORRELATION'': 0.9558373872353643}, ''times'': {''contamination'': 0.00020933151245117188, ''optimization'': 0.23014473915100098, ''imputation'': 0.04225325584411621, ''log_imputation'': -3.164073866748831}}, ''0.4'': {''scores'': {''RMSE'': 0.4508488005700101, ''MAE'': 0.34941433537269606, ''MI'': 0.8543113555966528, ''CORRELATION'': 0.8959297471926679}, ''times'': {''contamination'': 0.0003197193145751953, ''optimization'': 0.23014473915100098, ''imputation'': 0.046411752700805664, ''log_imputation'': -3.070202560857029}}, ''0.6'': {''scores'': {''RMSE'': 18.797539991079297, ''MAE'': 7.812583796335101, ''MI'': 0.36244773022350796, ''CORRELATION'': 0.6210142190959098}, ''times'': {''contamination'': 0.0006625652313232422, ''optimization'': 0.23014473915100098, ''imputation'': 0.043074607849121094, ''log_imputation'': -3.144821600518555}}, ''0.8'': {''scores'': {''RMSE'': 3.1451455567216193, ''MAE'': 1.1637520656636082, ''MI'': 0.0643204354315137, ''CORRELATION'': 0.22737088719870605}, ''times'': {''contamination'': 0.0013301372528076172, ''optimization'': 0.23014473915100098, ''imputation'': 0.03976869583129883, ''log_imputation'': -3.2246752130762237}}}}, ''iim'': {''bayesian'': {''0.05'': {''scores'': {''RMSE'': 0.2311363556202525, ''MAE'': 0.22809317150257158, ''MI'': 0.6931471805599452, ''CORRELATION'': 0.8754093900930757}, ''times'': {''contamination'': 0.00014829635620117188, ''optimization'': 0.8548746109008789, ''imputation'': 0.006244659423828125, ''log_imputation'': -5.0760286727078805}}, ''0.1'': {''scores'': {''RMSE'': 0.21734571962767568, ''MAE'': 0.20142183555276616, ''MI'': 1.4941751382893083, ''CORRELATION'': 0.9836625389334559}, ''times'': {''contamination'': 0.000431060791015625, ''optimization'': 0.8548746109008789, ''imputation'': 0.016748428344726562, ''log_imputation'': -4.089450855278965}}, ''0.2'': {''scores'': {''RMSE'': 0.2763681623559098, ''MAE'': 0.21205899863451294, ''MI'': 1.692828654044598, ''CORRELATION'': 0.9663556239228223}, ''times'': {''contamination'': 0.0001976490020751953, ''optimization'': 0.8548746109008789, ''imputation'': 0.038770198822021484, ''log_imputation'': -3.2501033990684074}}, ''0.4'': {''scores'': {''RMSE'': 0.32470532661816204, ''MAE'': 0.24836184775095202, ''MI'': 1.0631520030142667, ''CORRELATION'': 0.9435024215665483}, ''times'': {''contamination'': 0.0003495216369628906, ''optimization'': 0.8548746109008789, ''imputation'': 0.21004080772399902, ''log_imputation'': -1.5604534446476055}}, ''0.6'': {''scores'': {''RMSE'': 0.45693859713260937, ''MAE'': 0.3350566242376081, ''MI'': 0.836724518636222, ''CORRELATION'': 0.9015668975756113}, ''times'': {''contamination'': 0.0006783008575439453, ''optimization'': 0.8548746109008789, ''imputation'': 0.6758015155792236, ''log_imputation'': -0.39185586205985784}}, ''0.8'': {''scores'': {''RMSE'': 0.7301676007328138, ''MAE'': 0.5391664379693699, ''MI'': 0.43198783605819785, ''CORRELATION'': 0.7329833767632488}, ''times'': {''contamination'': 0.0021584033966064453, ''optimization'': 0.8548746109008789, ''imputation'': 3.9872961044311523, ''log_imputation'': 1.383113333118614}}}}, ''mrnn'': {''bayesian'': {''0.05'': {''scores'': {''RMSE'': 1.7800071295936744, ''MAE'': 1.7387530345617206, ''MI'': 0.0, ''CORRELATION'': 0}, ''times'': {''contamination'': 0.0002582073211669922, ''optimization'': 29.687320470809937, ''imputation'': 6.463064670562744, ''log_imputation'': 1.866103612506844}}, ''0.1'': {''scores'': {''RMSE'': 1.4143718875049651, ''MAE'': 1.2042932544309777, ''MI'': 0.0, ''CORRELATION'': 0}, ''times'': {''contamination'': 0.00035381317138671875, ''optimization'': 29.687320470809937, ''imputation'': 6.938860893249512, ''log_imputation'': 1.937137624627478}}, ''0.2'': {''scores'': {''RMSE'': 1.468628957888676, ''MAE'': 1.177928404642071, ''MI'': 0.2888172042752689, ''CORRELATION'': 0.27089827066191574}, ''times'': {''contamination'': 0.0003085136413574219, ''optimization'': 29.687320470809937, ''imputation'': 7.3980393409729, ''log_imputation'': 2.0012150111810514}}, ''0.4'': {''scores'': {''RMSE'': 1.2420980860848794, ''MAE'': 0.9793021593310115, ''MI'': 0.16817610173038405, ''CORRELATION'': 0.06862278658494796}, ''times'': {''contamination'': 0.0005781650543212891, ''optimization'': 29.687320470809937, ''imputation'': 7.899038076400757, ''log_imputation'': 2.066740989578406}}, ''0.6'': {''scores'': {''RMSE'': 1.3361811680109383, ''MAE'': 1.0908161906800318, ''MI'': 0.1307724816846375, ''CORRELATION'': 0.07966669775362903}, ''times'': {''contamination'': 0.0015287399291992188, ''optimization'': 29.687320470809937, ''imputation'': 5.665785551071167, ''log_imputation'': 1.7344455523106883}}, ''0.8'': {''scores'': {''RMSE'': 1.2352164167350344, ''MAE'': 0.9969833725134342, ''MI'': 0.08739712471151846, ''CORRELATION'': 0.03209263198570942}, ''times'': {''contamination'': 0.0018906593322753906, ''optimization'': 29.687320470809937, ''imputation'': 7.642791986465454, ''log_imputation'': 2.0337629796749472}}}}}}'

# --- RMSE sheet (mrnn rows, F34:F39) ---
$rmse = $wb.Worksheets.Item("RMSE")
$rmseRng = $rmse.Range("F34:F39")
$rmseStyle = $rmseRng.Style
$rmseRng.NumberFormat = "@"
$rmse.Range("F34").Value2 = "1.7800071296"
$rmse.Range("F35").Value2 = "1.4143718875"
$rmse.Range("F36").Value2 = "1.4686289579"
$rmse.Range("F37").Value2 = "1.2420980861"
$rmse.Range("F38").Value2 = "1.3361811680"
$rmse.Range("F39").Value2 = "1.2352164167"
$rmseRng.Style = $rmseStyle

# --- MAE sheet (mrnn rows, F34:F39) ---
$mae = $wb.Worksheets.Item("MAE")
$maeRng = $mae.Range("F34:F39")
$maeStyle = $maeRng.Style
$maeRng.NumberFormat = "@"
$mae.Range("F34").Value2 = "1.7387530346"
$mae.Range("F35").Value2 = "1.2042932544"
$mae.Range("F36").Value2 = "1.1779284046"
$mae.Range("F37").Value2 = "0.9793021593"
$mae.Range("F38").Value2 = "1.0908161907"
$mae.Range("F39").Value2 = "0.9969833725"
$maeRng.Style = $maeStyle

# --- MI sheet (mrnn rows, F36:F39) ---
$mi = $wb.Worksheets.Item("MI")
$miRng = $mi.Range("F36:F39")
$miStyle = $miRng.Style
$miRng.NumberFormat = "@"
$mi.Range("F36").Value2 = "0.2888172043"
$mi.Range("F37").Value2 = "0.1681761017"
$mi.Range("F38").Value2 = "0.1307724817"
$mi.Range("F39").Value2 = "0.0873971247"
$miRng.Style = $miStyle

# --- CORRELATION sheet (mrnn rows, F36:F39) ---
$corr = $wb.Worksheets.Item("CORRELATION")
$corrRng = $corr.Range("F36:F39")
$corrStyle = $corrRng.Style
$corrRng.NumberFormat = "@"
$corr.Range("F36").Value2 = "0.2708982707"
$corr.Range("F37").Value2 = "0.0686227866"
$corr.Range("F38").Value2 = "0.0796666978"
$corr.Range("F39").Value2 = "0.0320926320"
$corrRng.Style = $corrStyle
